$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.614.56"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "1.928.68"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'246.84"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.4750"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").Value = "'0.2924"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").Value = "'0.06818"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").Value = "'105.69"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").Value = "'18.47"
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").Value = "'0.07760"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "1.915.43"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'5.364"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "'0.6750"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'287.32"
$ws.Range("E16").Value = "  -6.56%  "
$ws.Range("D17").Value = "30.567.01"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.10"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007657"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "2.162.98"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'5.444"
$ws.Range("E22").Value = "  +3.95%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'6.286"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "'9.420"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'168.50"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'20.93"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'2.142"
$ws.Range("E28").Value = "  +7.12%  "
$ws.Range("D29").Value = "'0.1091"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'1.361"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'4.181"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'4.020"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'0.05080"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'0.7418"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "'1.156"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "'0.02095"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("D37").Value = "'2.730"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "'2.693"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").Value = "'2.070"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'111.16"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'0.8757"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").Value = "'0.4436"
$ws.Range("E42").Value = "  +4.57%  "
$ws.Range("D43").Value = "'5.919"
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'67.90"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "'7.299"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "'9.365"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "'47.98"
$ws.Range("E48").Value = "  +12.24%  "
$ws.Range("D49").Value = "'0.1240"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.4154"
$ws.Range("E50").Value = "  +7.11%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.32"
$ws.Range("E51").Value = "  +0.72%  "
